$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/artfynd/A 30234-2023.xlsx"; "A 30234-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/kartor/A 30234-2023.png; "A 30234-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomål/A 30234-2023.docx; "A 30234-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomålsmail/A 30234-2023.docx; "A 30234-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsyn/A 30234-2023.docx; "A 30234-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsynsmail/A 30234-2023.docx; "A 30234-2023")'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/artfynd/A 33548-2023.xlsx"; "A 33548-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/kartor/A 33548-2023.png; "A 33548-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomål/A 33548-2023.docx; "A 33548-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomålsmail/A 33548-2023.docx; "A 33548-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsyn/A 33548-2023.docx; "A 33548-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsynsmail/A 33548-2023.docx; "A 33548-2023")'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/artfynd/A 33550-2023.xlsx"; "A 33550-2023")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/kartor/A 33550-2023.png; "A 33550-2023")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomål/A 33550-2023.docx; "A 33550-2023")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomålsmail/A 33550-2023.docx; "A 33550-2023")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsyn/A 33550-2023.docx; "A 33550-2023")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsynsmail/A 33550-2023.docx; "A 33550-2023")'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/artfynd/A 30241-2023.xlsx"; "A 30241-2023")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/kartor/A 30241-2023.png; "A 30241-2023")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomål/A 30241-2023.docx; "A 30241-2023")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomålsmail/A 30241-2023.docx; "A 30241-2023")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsyn/A 30241-2023.docx; "A 30241-2023")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsynsmail/A 30241-2023.docx; "A 30241-2023")'
